$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old detail rows (3-6); only the summary row (row 2) remains.
$ws.Range("A3:H6").Clear()

# Drop the "Unnamed: 0.1" index column (old column B) and shift everything
# else (old C:H) one column to the left.
$ws.Range("B1").EntireColumn.Delete()

# Fix up the remaining data row with the corrected transaction values.
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Direct Deposit"
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "'1/13/2025"
$ws.Range("G2").Style = $ws.Range("F2").Style
